$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the four new "correlation threshold" sheets at the end of the
#    workbook (after "Уровень значимости").
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Корреляция более 75%"

$ws4 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Корреляция более 80%"

$ws5 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "Корреляция более 85%"

$ws6 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws6.Name = "Корреляция более 90%"

# ---------------------------------------------------------------------------
# 2. Prime the shared-string table on sheet3 so the new strings land at the
#    exact indices used elsewhere (820-835): all "forward" pair descriptions
#    first (a<b), then all "reversed" pair descriptions (a>b) - both in the
#    sheet's final row order. The scratch column is cleared right after.
# ---------------------------------------------------------------------------
$forward = @(
  "3 DER_pt_h 9 DER_sum_pt (0.80465631268380367, 0.0)",
  "3 DER_pt_h 21 PRI_met_sumet (0.76295643297226268, 0.0)",
  "3 DER_pt_h 29 PRI_jet_all_pt (0.78030514800037387, 0.0)",
  "9 DER_sum_pt 21 PRI_met_sumet (0.90136734493390336, 0.0)",
  "9 DER_sum_pt 22 PRI_jet_num (0.77381750235404123, 0.0)",
  "9 DER_sum_pt 29 PRI_jet_all_pt (0.96300530649741123, 0.0)",
  "21 PRI_met_sumet 29 PRI_jet_all_pt (0.88153096941468034, 0.0)",
  "22 PRI_jet_num 29 PRI_jet_all_pt (0.81379266838833642, 0.0)"
)
$backward = @(
  "9 DER_sum_pt 3 DER_pt_h (0.80465631268380367, 0.0)",
  "21 PRI_met_sumet 3 DER_pt_h (0.76295643297226268, 0.0)",
  "21 PRI_met_sumet 9 DER_sum_pt (0.90136734493390336, 0.0)",
  "22 PRI_jet_num 9 DER_sum_pt (0.77381750235404123, 0.0)",
  "29 PRI_jet_all_pt 3 DER_pt_h (0.78030514800037387, 0.0)",
  "29 PRI_jet_all_pt 9 DER_sum_pt (0.96300530649741123, 0.0)",
  "29 PRI_jet_all_pt 21 PRI_met_sumet (0.88153096941468034, 0.0)",
  "29 PRI_jet_all_pt 22 PRI_jet_num (0.81379266838833642, 0.0)"
)

$primeRow = 1
foreach ($s in $forward) {
  $ws3.Cells.Item($primeRow, 26).Value = $s
  $primeRow = $primeRow + 1
}
foreach ($s in $backward) {
  $ws3.Cells.Item($primeRow, 26).Value = $s
  $primeRow = $primeRow + 1
}
$ws3.Columns.Item(26).ClearContents()

# ---------------------------------------------------------------------------
# 3. Fill in the real content of "Корреляция более 75%" (A1:M16).
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "3 DER_pt_h 9 DER_sum_pt (0.80465631268380367, 0.0)"
$ws3.Range("J1").Value = 3
$ws3.Range("K1").Value = 9
$ws3.Range("L1").Value = 21
$ws3.Range("M1").Value = 29

$ws3.Range("A2").Value = "3 DER_pt_h 21 PRI_met_sumet (0.76295643297226268, 0.0)"
$ws3.Range("J2").Value = 9
$ws3.Range("K2").Value = 22

$ws3.Range("A3").Value = "3 DER_pt_h 29 PRI_jet_all_pt (0.78030514800037387, 0.0)"
$ws3.Range("J3").Value = 29
$ws3.Range("K3").Value = 22

$ws3.Range("A4").Value = "9 DER_sum_pt 3 DER_pt_h (0.80465631268380367, 0.0)"
$ws3.Range("A5").Value = "9 DER_sum_pt 21 PRI_met_sumet (0.90136734493390336, 0.0)"
$ws3.Range("A6").Value = "9 DER_sum_pt 22 PRI_jet_num (0.77381750235404123, 0.0)"
$ws3.Range("A7").Value = "9 DER_sum_pt 29 PRI_jet_all_pt (0.96300530649741123, 0.0)"
$ws3.Range("A8").Value = "21 PRI_met_sumet 3 DER_pt_h (0.76295643297226268, 0.0)"
$ws3.Range("A9").Value = "21 PRI_met_sumet 9 DER_sum_pt (0.90136734493390336, 0.0)"
$ws3.Range("A10").Value = "21 PRI_met_sumet 29 PRI_jet_all_pt (0.88153096941468034, 0.0)"
$ws3.Range("A11").Value = "22 PRI_jet_num 9 DER_sum_pt (0.77381750235404123, 0.0)"
$ws3.Range("A12").Value = "22 PRI_jet_num 29 PRI_jet_all_pt (0.81379266838833642, 0.0)"
$ws3.Range("A13").Value = "29 PRI_jet_all_pt 3 DER_pt_h (0.78030514800037387, 0.0)"
$ws3.Range("A14").Value = "29 PRI_jet_all_pt 9 DER_sum_pt (0.96300530649741123, 0.0)"
$ws3.Range("A15").Value = "29 PRI_jet_all_pt 21 PRI_met_sumet (0.88153096941468034, 0.0)"
$ws3.Range("A16").Value = "29 PRI_jet_all_pt 22 PRI_jet_num (0.81379266838833642, 0.0)"

# ---------------------------------------------------------------------------
# 4. Fill in "Корреляция более 80%" (A1:L10).
# ---------------------------------------------------------------------------
$ws4.Range("A1").Value = "3 DER_pt_h 9 DER_sum_pt (0.80465631268380367, 0.0)"
$ws4.Range("J1").Value = 3
$ws4.Range("K1").Value = 9

$ws4.Range("A2").Value = "9 DER_sum_pt 3 DER_pt_h (0.80465631268380367, 0.0)"
$ws4.Range("J2").Value = 9
$ws4.Range("K2").Value = 21
$ws4.Range("L2").Value = 29

$ws4.Range("A3").Value = "9 DER_sum_pt 21 PRI_met_sumet (0.90136734493390336, 0.0)"
$ws4.Range("J3").Value = 22
$ws4.Range("K3").Value = 29

$ws4.Range("A4").Value = "9 DER_sum_pt 29 PRI_jet_all_pt (0.96300530649741123, 0.0)"
$ws4.Range("A5").Value = "21 PRI_met_sumet 9 DER_sum_pt (0.90136734493390336, 0.0)"
$ws4.Range("A6").Value = "21 PRI_met_sumet 29 PRI_jet_all_pt (0.88153096941468034, 0.0)"
$ws4.Range("A7").Value = "22 PRI_jet_num 29 PRI_jet_all_pt (0.81379266838833642, 0.0)"
$ws4.Range("A8").Value = "29 PRI_jet_all_pt 9 DER_sum_pt (0.96300530649741123, 0.0)"
$ws4.Range("A9").Value = "29 PRI_jet_all_pt 21 PRI_met_sumet (0.88153096941468034, 0.0)"
$ws4.Range("A10").Value = "29 PRI_jet_all_pt 22 PRI_jet_num (0.81379266838833642, 0.0)"

# ---------------------------------------------------------------------------
# 5. Fill in "Корреляция более 85%" (A1:L6).
# ---------------------------------------------------------------------------
$ws5.Range("A1").Value = "9 DER_sum_pt 21 PRI_met_sumet (0.90136734493390336, 0.0)"
$ws5.Range("J1").Value = 9
$ws5.Range("K1").Value = 21
$ws5.Range("L1").Value = 29

$ws5.Range("A2").Value = "9 DER_sum_pt 29 PRI_jet_all_pt (0.96300530649741123, 0.0)"
$ws5.Range("A3").Value = "21 PRI_met_sumet 9 DER_sum_pt (0.90136734493390336, 0.0)"
$ws5.Range("A4").Value = "21 PRI_met_sumet 29 PRI_jet_all_pt (0.88153096941468034, 0.0)"
$ws5.Range("A5").Value = "29 PRI_jet_all_pt 9 DER_sum_pt (0.96300530649741123, 0.0)"
$ws5.Range("A6").Value = "29 PRI_jet_all_pt 21 PRI_met_sumet (0.88153096941468034, 0.0)"

# ---------------------------------------------------------------------------
# 6. Fill in "Корреляция более 90%" (A1:L4).
# ---------------------------------------------------------------------------
$ws6.Range("A1").Value = "9 DER_sum_pt 21 PRI_met_sumet (0.90136734493390336, 0.0)"
$ws6.Range("J1").Value = 9
$ws6.Range("K1").Value = 21
$ws6.Range("L1").Value = 29

$ws6.Range("A2").Value = "9 DER_sum_pt 29 PRI_jet_all_pt (0.96300530649741123, 0.0)"
$ws6.Range("A3").Value = "21 PRI_met_sumet 9 DER_sum_pt (0.90136734493390336, 0.0)"
$ws6.Range("A4").Value = "29 PRI_jet_all_pt 9 DER_sum_pt (0.96300530649741123, 0.0)"

# ---------------------------------------------------------------------------
# 7. Selections / active-tab bookkeeping: activate sheet1 and sheet2 first
#    (updating their selections without leaving tabSelected set on them),
#    then finally activate sheet6 so it becomes the single active tab.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A6:XFD6").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("B1").Select()

$ws6.Activate()

# Best-effort: scroll the tab strip so the first visible tab is index 1
# (matches firstSheet="1" in the saved bookView).
try { $wb.Windows.Item(1).ScrollWorkbookTabs(1) | Out-Null } catch {}
try { $wb.Windows.Item(1).DisplayedTab = 1 } catch {}
